$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 / J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style of the existing header row (copy from H1, which already
# carries style index 1: bold font, thin box border, centered/top alignment).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data for the new I/J columns (rows 2-36)
$values = @(
  @(8, 8),
  @(5, 5),
  @(6, 7),
  @(8, 8),
  @(8, 8),
  @(6, 6),
  @(9, 9),
  @(7, 8),
  @(8, 8),
  @(4, 4),
  @(8, 8),
  @(6, 6),
  @(8, 8),
  @(5, 6),
  @(5, 5),
  @(10, 10),
  @(10, 10),
  @(6, 6),
  @(7, 8),
  @(7, 7),
  @(7, 7),
  @(5, 5),
  @(9, 9),
  @(6, 6),
  @(6, 8),
  @(5, 5),
  @(8, 8),
  @(7, 7),
  @(4, 5),
  @(6, 6),
  @(5, 5),
  @(5, 5),
  @(1, 1),
  @(5, 5),
  @(4, 4)
)

for ($i = 0; $i -lt $values.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 9).Value = $values[$i][0]
  $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
